$wb = $excel.ActiveWorkbook

# Update the "展览" (Exhibition) sheet
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1030
$ws1.Range("F3").Value = 15
$ws1.Range("F4").Value = 502

# Update the "全部类型" (All types) sheet
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1030
$ws4.Range("F3").Value = 15
$ws4.Range("F4").Value = 502
